$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: row, A, B, D, E, F, G, H, K, L, M, N, Q(date)
$data = @(
    @(2, 34, 6475, 4, 55043, 1583242, 21325, 13400, 11678, 316002, 6011, 3121, "2018-08-14 12:17:01"),
    @(3, 35, 7755, 5, 66726, 1899239, 27596, 16677, 11683, 315997, 6271, 3277, "2018-08-14 12:17:11"),
    @(4, 36, 9035, 6, 74051, 2219594, 30505, 19141, 7325, 320355, 2909, 2464, "2018-08-14 12:17:21"),
    @(5, 37, 10315, 7, 85984, 2535341, 37160, 22464, 11933, 315747, 6655, 3323, "2018-08-14 12:17:31"),
    @(6, 38, 11595, 8, 93340, 2855665, 40068, 24720, 7356, 320324, 2908, 2256, "2018-08-14 12:17:41"),
    @(7, 39, 12875, 9, 105055, 3171631, 46657, 28326, 11715, 315966, 6589, 3606, "2018-08-14 12:17:51"),
    @(8, 40, 14155, 10, 112410, 3491955, 50017, 30987, 7355, 320324, 3360, 2661, "2018-08-14 12:18:01"),
    @(9, 41, 15435, 11, 124264, 3807781, 56281, 34477, 11854, 315826, 6264, 3490, "2018-08-14 12:18:11"),
    @(10, 42, 16715, 12, 136019, 4123685, 62545, 37786, 11754, 315904, 6264, 3309, "2018-08-14 12:18:21"),
    @(11, 43, 17995, 13, 143448, 4443957, 65839, 40405, 7428, 320272, 3294, 2619, "2018-08-14 12:18:31"),
    @(12, 44, 19275, 14, 155104, 4759981, 72107, 43837, 11656, 316024, 6268, 3432, "2018-08-14 12:18:41"),
    @(13, 45, 20555, 15, 166842, 5075902, 78571, 47121, 11738, 315921, 6464, 3284, "2018-08-14 12:18:51"),
    @(14, 46, 21835, 16, 174237, 5396208, 81993, 49833, 7394, 320306, 3422, 2712, "2018-08-14 12:19:01"),
    @(15, 47, 23115, 17, 186033, 5712092, 87810, 52931, 11796, 315884, 5817, 3098, "2018-08-14 12:19:11"),
    @(16, 48, 24395, 18, 193390, 6032415, 90718, 55408, 7357, 320323, 2908, 2477, "2018-08-14 12:19:21"),
    @(17, 49, 25675, 19, 205039, 6348446, 96534, 58614, 11649, 316031, 5816, 3206, "2018-08-14 12:19:31"),
    @(18, 50, 26955, 20, 216696, 6664469, 102352, 61659, 11657, 316023, 5818, 3045, "2018-08-14 12:19:41"),
    @(19, 51, 28235, 21, 224057, 6984788, 105456, 64219, 7361, 320319, 3104, 2560, "2018-08-14 12:19:51"),
    @(20, 52, 29515, 22, 235791, 7300734, 111398, 67216, 11733, 315946, 5942, 2997, "2018-08-14 12:20:01"),
    @(21, 53, 30795, 23, 247462, 7616743, 117214, 70258, 11671, 316009, 5816, 3042, "2018-08-14 12:20:11"),
    @(22, 54, 32075, 24, 254838, 7937047, 120508, 72776, 7376, 320304, 3294, 2518, "2018-08-14 12:20:21"),
    @(23, 55, 33355, 25, 266580, 8252985, 126452, 75785, 11742, 315938, 5944, 3009, "2018-08-14 12:20:31"),
    @(24, 56, 34635, 26, 278310, 8568914, 132982, 79303, 11730, 315929, 6530, 3518, "2018-08-14 12:20:41"),
    @(25, 57, 35915, 27, 285708, 8889217, 136080, 81737, 7397, 320303, 3098, 2434, "2018-08-14 12:20:51"),
    @(26, 58, 37195, 28, 297417, 9205188, 142220, 84751, 11709, 315971, 6140, 3014, "2018-08-14 12:21:01"),
    @(27, 59, 38475, 29, 304793, 9525492, 145192, 87134, 7375, 320304, 2972, 2383, "2018-08-14 12:21:11"),
    @(28, 60, 39755, 30, 316500, 9841465, 151394, 90368, 11707, 315973, 6202, 3234, "2018-08-14 12:21:21"),
    @(29, 61, 41035, 31, 328220, 10157404, 157917, 93936, 11720, 315939, 6523, 3568, "2018-08-14 12:21:31"),
    @(30, 62, 42315, 32, 335706, 10477619, 160826, 96241, 7485, 320215, 2909, 2305, "2018-08-14 12:21:41"),
    @(31, 63, 43595, 33, 347545, 10793460, 166708, 99184, 11839, 315841, 5882, 2943, "2018-08-14 12:21:51"),
    @(32, 64, 44875, 34, 354947, 11113738, 169679, 101576, 7402, 320278, 2971, 2392, "2018-08-14 12:22:01"),
    @(33, 65, 46155, 35, 366695, 11429649, 176074, 104777, 11748, 315911, 6395, 3201, "2018-08-14 12:22:11"),
    @(34, 66, 47435, 36, 374067, 11749978, 178982, 107143, 7371, 320329, 2908, 2366, "2018-08-14 12:22:21")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $row[1]    # A  id
    $ws.Cells.Item($r, 2).Value  = $row[2]    # B  clock_time
    $ws.Cells.Item($r, 4).Value  = $row[3]    # D  seqno
    $ws.Cells.Item($r, 5).Value  = $row[4]    # E  all_cpu
    $ws.Cells.Item($r, 6).Value  = $row[5]    # F  all_lpm
    $ws.Cells.Item($r, 7).Value  = $row[6]    # G  all_transmit
    $ws.Cells.Item($r, 8).Value  = $row[7]    # H  all_listen
    $ws.Cells.Item($r, 11).Value = $row[8]    # K  cpu
    $ws.Cells.Item($r, 12).Value = $row[9]    # L  lpm
    $ws.Cells.Item($r, 13).Value = $row[10]   # M  transmit
    $ws.Cells.Item($r, 14).Value = $row[11]   # N  listen
    $ws.Cells.Item($r, 17).Value = $row[12]   # Q  date
}
